# Update "想去人数" (want-to-go count) figures and one cover-image URL
# to match the newly generated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

function Update-SheetCounts {
    param(
        [string]$SheetName,
        [hashtable]$FUpdates,
        [hashtable]$IUpdates
    )

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in $FUpdates.Keys) {
        $ws.Range("F$row").Value = $FUpdates[$row]
    }

    foreach ($row in $IUpdates.Keys) {
        $ws.Range("I$row").Value = $IUpdates[$row]
    }
}

# Sheet "展览" (exhibitions)
$exhibitionF = @{
    5  = 1817
    6  = 460
    8  = 163
    9  = 2346
    10 = 126
    11 = 68
    13 = 1420
    14 = 504
    17 = 220
    21 = 209
    24 = 88
    25 = 32
    26 = 1459
    28 = 369
    29 = 219
    31 = 286
    32 = 364
}
$exhibitionI = @{
    29 = "//i0.hdslb.com/bfs/openplatform/202405/fgnndv151716373477064.jpeg"
}
Update-SheetCounts "展览" $exhibitionF $exhibitionI

# Sheet "全部类型" (all types)
$allTypesF = @{
    5  = 1817
    7  = 460
    9  = 163
    10 = 2346
    11 = 126
    12 = 68
    14 = 1420
    15 = 504
    18 = 220
    22 = 209
    25 = 88
    26 = 32
    27 = 1459
    29 = 369
    30 = 219
    32 = 286
    33 = 364
}
$allTypesI = @{
    30 = "//i0.hdslb.com/bfs/openplatform/202405/fgnndv151716373477064.jpeg"
}
Update-SheetCounts "全部类型" $allTypesF $allTypesI
